# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Most "Price" (column D) values are plain decimals that Excel would otherwise
# auto-convert to numbers on assignment, so we force those cells to Text
# (NumberFormat "@") before writing, keeping them as literal strings exactly
# like the original inlineStr cells. Values that already contain a second
# "." (thousands separator, e.g. "36.935.50") or other non-numeric
# characters are left alone since Excel can't parse them as a number anyway.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.935.50"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.041.77"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.76"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.657"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.78"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0767"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.44"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.877"
$ws.Range("E13").Value = "  +9.99%  "
$ws.Range("D14").Value = "2.341.93"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.63"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "2.030.04"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.05"
$ws.Range("E17").Value = "  +9.49%  "
$ws.Range("D18").Value = "36.931.16"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.68"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.38"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.42"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +3.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.61"
$ws.Range("E25").Value = "  +5.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.38"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.12"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.87"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.33"
$ws.Range("E29").Value = "  +14.44%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.67"
$ws.Range("E32").Value = "  +5.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0611"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0871"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").Value = "  +6.16%  "
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.07"
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0982"
$ws.Range("E41").Value = "  -7.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0221"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.13"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.62"
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.81"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").Value = "1.289.56"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.78"
$ws.Range("E47").Value = "  +8.53%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.71"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").Value = "2.226.53"
$ws.Range("E51").Value = "  -0.22%  "
